$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.279.57'
$ws.Range("E2").Value = '  -1.02%  '

# Row 3
$ws.Range("D3").Value = '1.705.73'
$ws.Range("E3").Value = '  -1.07%  '

# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '''223.96'
$ws.Range("E5").Value = '  -1.01%  '

# Row 6
$ws.Range("D6").Value = '''0.5318'
$ws.Range("E6").Value = '  -1.14%  '

# Row 7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").Value = '''0.2661'
$ws.Range("E8").Value = '  -0.52%  '

# Row 9
$ws.Range("D9").Value = '''0.06600'
$ws.Range("E9").Value = '  -0.20%  '

# Row 10
$ws.Range("D10").Value = '''20.73'
$ws.Range("E10").Value = '  -4.86%  '

# Row 11
$ws.Range("D11").Value = '''0.07674'
$ws.Range("E11").Value = '  -0.61%  '

# Row 12
$ws.Range("D12").Value = '''4.506'
$ws.Range("E12").Value = '  -2.34%  '

# Row 13
$ws.Range("D13").Value = '1.714.53'
$ws.Range("E13").Value = '  -0.60%  '

# Row 14
$ws.Range("D14").Value = '1.940.14'
$ws.Range("E14").Value = '  -1.13%  '

# Row 15
$ws.Range("D15").Value = '''0.5815'
$ws.Range("E15").Value = '  -1.03%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8169'
$ws.Range("E16").Value = '  -1.82%  '

# Row 17
$ws.Range("D17").Value = '''67.59'
$ws.Range("E17").Value = '  -0.68%  '

# Row 18
$ws.Range("D18").Value = '27.296.60'
$ws.Range("E18").Value = '  -1.03%  '

# Row 19
$ws.Range("D19").Value = '''215.34'
$ws.Range("E19").Value = '  -2.99%  '

# Row 20
$ws.Range("E20").Value = '  -0.05%  '

# Row 21
$ws.Range("D21").Value = '''4.625'
$ws.Range("E21").Value = '  -2.45%  '

# Row 22
$ws.Range("D22").Value = '''10.41'
$ws.Range("E22").Value = '  -2.72%  '

# Row 23
$ws.Range("D23").Value = '''5.982'
$ws.Range("E23").Value = '  -2.01%  '

# Row 24
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").Value = '''143.91'
$ws.Range("E25").Value = '  -2.90%  '

# Row 26
$ws.Range("D26").Value = '''1.690'
$ws.Range("E26").Value = '  -0.31%  '

# Row 27
$ws.Range("D27").Value = '''0.1205'
$ws.Range("E27").Value = '  -2.47%  '

# Row 28
$ws.Range("D28").Value = '''7.224'
$ws.Range("E28").Value = '  -2.42%  '

# Row 29
$ws.Range("D29").Value = '''16.19'
$ws.Range("E29").Value = '  -2.87%  '

# Row 30
$ws.Range("D30").Value = '''0.05364'
$ws.Range("E30").Value = '  -3.33%  '

# Row 31
$ws.Range("D31").Value = '''1.289'
$ws.Range("E31").Value = '  -1.17%  '

# Row 32
$ws.Range("D32").Value = '''3.484'
$ws.Range("E32").Value = '  -1.78%  '

# Row 33
$ws.Range("D33").Value = '''3.416'
$ws.Range("E33").Value = '  -1.31%  '

# Row 34
$ws.Range("D34").Value = '''1.643'
$ws.Range("E34").Value = '  -1.23%  '

# Row 35
$ws.Range("E35").Value = '  +1.48%  '

# Row 36
$ws.Range("D36").Value = '''0.9502'
$ws.Range("E36").Value = '  -1.39%  '

# Row 37
$ws.Range("D37").Value = '''2.398'
$ws.Range("E37").Value = '  -1.89%  '

# Row 38
$ws.Range("D38").Value = '''0.5842'
$ws.Range("E38").Value = '  -2.01%  '

# Row 39
$ws.Range("D39").Value = '''0.01640'
$ws.Range("E39").Value = '  -0.51%  '

# Row 40
$ws.Range("D40").Value = '''5.810'
$ws.Range("E40").Value = '  -2.01%  '

# Row 41
$ws.Range("D41").Value = '1.043.52'
$ws.Range("E41").Value = '  -1.36%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.003'
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.8420'
$ws.Range("E43").Value = '  -1.60%  '

# Row 44
$ws.Range("D44").Value = '''100.89'
$ws.Range("E44").Value = '  -0.50%  '

# Row 45
$ws.Range("D45").Value = '1.848.14'
$ws.Range("E45").Value = '  -1.09%  '

# Row 46
$ws.Range("E46").Value = '  -4.29%  '

# Row 47
$ws.Range("D47").Value = '''57.84'
$ws.Range("E47").Value = '  -2.21%  '

# Row 48
$ws.Range("D48").Value = '''0.4522'
$ws.Range("E48").Value = '  +1.82%  '

# Row 49
$ws.Range("D49").Value = '''1.005'
$ws.Range("E49").Value = '  +0.11%  '

# Row 50
$ws.Range("D50").Value = '''8.101'
$ws.Range("E50").Value = '  -1.50%  '

# Row 51
$ws.Range("D51").Value = '''0.05235'
$ws.Range("E51").Value = '  -0.72%  '

